# Sciences.xlsx - Spring 2026 (1261) term requisites refresh.
# The course list gained 4 new rows (SCI100, SCI130, SCI140, SCI834) and the
# whole table was re-sorted alphabetically by Course_Code, so the simplest
# faithful reproduction is: clear the old data body and rewrite it complete,
# in the new order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (rows 2-7) before rewriting with the new, re-sorted,
# expanded set (rows 2-11).
$ws.Range("A2:I7").ClearContents()

$data = @(
  @("SCI100","Course","Sciences","SCI 100 - Success Strategies in Science","None","None","SCI","DEAN SCI",""),
  @("SCI130","Course","Sciences","SCI 130 - Busting Science Myths: Critical Thinking in the Misinformation Age","None","None","SCI","DEAN SCI",""),
  @("SCI140","Course","Sciences","SCI 140 - Genes Unzipped: The Science and Secrets of DNA","None","None","SCI","DEAN SCI",""),
  @("SCI190","Course","Sciences","SCI 190 - Special Topics in Science Interdisciplinarity","None","None","SCI","DEAN SCI",""),
  @("SCI191","Course","Sciences","SCI 191 - Introduction to Modern Scientific Research","None","None","SCI","DEAN SCI",""),
  @("SCI192","Course","Sciences","SCI 192 - 'The Science Around Us' Interdisciplinary Science for Pre-Service Teachers","None","None","SCI","DEAN SCI",""),
  @("SCI300","Course","Sciences","SCI 300 - Science and its Impact on Society","None","None","SCI","DEAN SCI","REQ-60 units.  Not open to students in the Faculty of Science or the Schools of Computing Science, Engineering Science."),
  @("SCI301","Course","Sciences","SCI 301 - Science Communication: An Introduction","None","None","SCI","DEAN SCI","REQ-60 units towards a BSc degree or permission of instructor."),
  @("SCI390","Course","Sciences","SCI 390 - Sepcial Topics in Science Interdisciplinarity","None","None","SCI","DEAN SCI",""),
  @("SCI834","Course","Sciences","SCI 834 - Essential Cell Biology","None","None","SCI","DEAN SCI","")
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 9; $c++) {
        $val = $row[$c - 1]
        # Description (column I) is blank for most courses - leave those
        # cells untouched/empty rather than writing an empty string, so the
        # sheet doesn't grow a spurious shared-string cell there.
        if ($val -ne "") {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
    $r++
}

# Matches the author's final selection after entering the new data.
$ws.Range("A2:I11").Select()
